# "Add files via upload" -- the CPU and GPU sheets each receive additional
# GPU benchmark rows that already exist (as data) on the SPEC sheet, and the
# active-sheet / selection state left behind by the edit changes too.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# CPU sheet (sheet1): append rows 4-6
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("CPU")

$ws1.Cells.Item(4,1).Value = "RX 7801 XT"
$ws1.Cells.Item(4,2).Value = 111
$ws1.Cells.Item(4,3).Value = 91
$ws1.Cells.Item(4,4).Value = 56

$ws1.Cells.Item(5,1).Value = "RX 7802 XT"
$ws1.Cells.Item(5,2).Value = 112
$ws1.Cells.Item(5,3).Value = 92
$ws1.Cells.Item(5,4).Value = 57

$ws1.Cells.Item(6,1).Value = "RX 7813 XT"
$ws1.Cells.Item(6,2).Value = 123
$ws1.Cells.Item(6,3).Value = 103
$ws1.Cells.Item(6,4).Value = 68

# ---------------------------------------------------------------------
# GPU sheet (sheet2): append rows 4-16
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("GPU")

$ws2.Cells.Item(4,1).Value = "RX 7801 XT"
$ws2.Cells.Item(4,2).Value = 111
$ws2.Cells.Item(4,3).Value = 91
$ws2.Cells.Item(4,4).Value = 56

$ws2.Cells.Item(5,1).Value = "RX 7802 XT"
$ws2.Cells.Item(5,2).Value = 112
$ws2.Cells.Item(5,3).Value = 92
$ws2.Cells.Item(5,4).Value = 57

$ws2.Cells.Item(6,1).Value = "RX 7803 XT"
$ws2.Cells.Item(6,2).Value = 113
$ws2.Cells.Item(6,3).Value = 93
$ws2.Cells.Item(6,4).Value = 58

$ws2.Cells.Item(7,1).Value = "RX 7804 XT"
$ws2.Cells.Item(7,2).Value = 114
$ws2.Cells.Item(7,3).Value = 94
$ws2.Cells.Item(7,4).Value = 59

$ws2.Cells.Item(8,1).Value = "RX 7805 XT"
$ws2.Cells.Item(8,2).Value = 115
$ws2.Cells.Item(8,3).Value = 95
$ws2.Cells.Item(8,4).Value = 60

$ws2.Cells.Item(9,1).Value = "RX 7806 XT"
$ws2.Cells.Item(9,2).Value = 116
$ws2.Cells.Item(9,3).Value = 96
$ws2.Cells.Item(9,4).Value = 61

$ws2.Cells.Item(10,1).Value = "RX 7807 XT"
$ws2.Cells.Item(10,2).Value = 117
$ws2.Cells.Item(10,3).Value = 97
$ws2.Cells.Item(10,4).Value = 62

$ws2.Cells.Item(11,1).Value = "RX 7808 XT"
$ws2.Cells.Item(11,2).Value = 118
$ws2.Cells.Item(11,3).Value = 98
$ws2.Cells.Item(11,4).Value = 63

$ws2.Cells.Item(12,1).Value = "RX 7809 XT"
$ws2.Cells.Item(12,2).Value = 119
$ws2.Cells.Item(12,3).Value = 99
$ws2.Cells.Item(12,4).Value = 64

$ws2.Cells.Item(13,1).Value = "RX 7810 XT"
$ws2.Cells.Item(13,2).Value = 120
$ws2.Cells.Item(13,3).Value = 100
$ws2.Cells.Item(13,4).Value = 65

$ws2.Cells.Item(14,1).Value = "RX 7811 XT"
$ws2.Cells.Item(14,2).Value = 121
$ws2.Cells.Item(14,3).Value = 101
$ws2.Cells.Item(14,4).Value = 66

$ws2.Cells.Item(15,1).Value = "RX 7812 XT"
$ws2.Cells.Item(15,2).Value = 122
$ws2.Cells.Item(15,3).Value = 102
$ws2.Cells.Item(15,4).Value = 67

$ws2.Cells.Item(16,1).Value = "RX 7813 XT"
$ws2.Cells.Item(16,2).Value = 123
$ws2.Cells.Item(16,3).Value = 103
$ws2.Cells.Item(16,4).Value = 68

# ---------------------------------------------------------------------
# Leave behind the same view/selection state as the saved workbook:
#   - GPU:  cell F11 selected (not the active tab)
#   - SPEC: A1:D16 selected, not the active tab anymore
#   - CPU:  D8 selected and CPU is the active tab on reopen
# (selecting a range activates its sheet/tab, so GPU and SPEC are
# touched first and CPU is activated last so it "wins" as active tab.)
# ---------------------------------------------------------------------
$ws2.Range("F11").Select()

$ws4 = $wb.Worksheets.Item("SPEC")
$ws4.Range("A1:D16").Select()
$ws4.Range("D16").Select()

$ws1.Range("D8").Select()
$ws1.Activate()
